$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A19").Value = 42624.619699074072
$ws.Range("A19").NumberFormat = "m/d/yy h:mm"

$ws.Range("B19").Value = 30
$ws.Range("C19").Value = 55
$ws.Range("D19").Value = 44
$ws.Range("E19").Value = 55
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = 8497
$ws.Range("H19").Value = 5514
$ws.Range("I19").Value = 1072
$ws.Range("J19").Value = 116
$ws.Range("K19").Value = 92
$ws.Range("L19").Value = 32
$ws.Range("M19").Value = 3
$ws.Range("N19").Value = "Noun"

$wb.Save()
